$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1289.8889
$ws.Range("I41").Value = 2128.4285
$ws.Range("J41").Value = 756.2727
$ws.Range("K41").Value = 2128.4285
$ws.Range("L41").Value = 756.2727
$ws.Range("M41").Value = -1688.4285
$ws.Range("N41").Value = -1636.2727
$ws.Range("H118").Value = 7909.2
$ws.Range("I118").Value = 11334.8
$ws.Range("J118").Value = 1058
$ws.Range("K118").Value = 34004.39999999999
$ws.Range("L118").Value = 3174
$ws.Range("M118").Value = -32347.39999999999
$ws.Range("N118").Value = -6488
$ws.Range("H127").Value = 21278806
$ws.Range("I127").Value = 437.5
$ws.Range("J127").Value = 24392714
$ws.Range("K127").Value = 1312.5
$ws.Range("L127").Value = 73178142
$ws.Range("M127").Value = 3647.5
$ws.Range("N127").Value = -73188062
$ws.Range("H129").Value = 899.9524
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 899.9524
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 2699.8572
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -12699.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25431.846
$ws.Range("I32").Value = 4327.952
$ws.Range("K32").Value = 4327.952
$ws.Range("M32").Value = -4040.952
$ws.Range("H61").Value = 1064.5186
$ws.Range("I61").Value = 963.7083
$ws.Range("J61").Value = 1871
$ws.Range("K61").Value = 963.7083
$ws.Range("L61").Value = 1871
$ws.Range("M61").Value = -751.7083
$ws.Range("N61").Value = -2295
$ws.Range("H88").Value = 2861.2
$ws.Range("I88").Value = 2944.5715
$ws.Range("J88").Value = 2666.6667
$ws.Range("K88").Value = 2944.5715
$ws.Range("L88").Value = 2666.6667
$ws.Range("M88").Value = -2538.5715
$ws.Range("N88").Value = -3478.6667
$ws.Range("H91").Value = 2861.2
$ws.Range("I91").Value = 2944.5715
$ws.Range("J91").Value = 2666.6667
$ws.Range("K91").Value = 2944.5715
$ws.Range("L91").Value = 2666.6667
$ws.Range("M91").Value = -1540.5715
$ws.Range("N91").Value = -5474.6667
$ws.Range("H136").Value = 1064.5186
$ws.Range("I136").Value = 963.7083
$ws.Range("J136").Value = 1871
$ws.Range("K136").Value = 2891.1249
$ws.Range("L136").Value = 5613
$ws.Range("M136").Value = -341.1248999999998
$ws.Range("N136").Value = -10713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 59528.05
$ws.Range("I86").Value = 86036.38
$ws.Range("J86").Value = 2093.3333
$ws.Range("K86").Value = 86036.38
$ws.Range("L86").Value = 2093.3333
$ws.Range("M86").Value = -84913.38
$ws.Range("N86").Value = -4339.3333
$ws.Range("H89").Value = 59528.05
$ws.Range("I89").Value = 86036.38
$ws.Range("J89").Value = 2093.3333
$ws.Range("K89").Value = 430181.9
$ws.Range("L89").Value = 10466.6665
$ws.Range("M89").Value = -424565.9
$ws.Range("N89").Value = -21698.6665
$ws.Range("H107").Value = 55556264
$ws.Range("I107").Value = 76923710
$ws.Range("K107").Value = 76923710
$ws.Range("M107").Value = -76921790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 998.3333
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -713
$ws.Range("H31").Value = 37128.723
$ws.Range("I31").Value = 68282.2
$ws.Range("J31").Value = 3750
$ws.Range("K31").Value = 68282.2
$ws.Range("L31").Value = 3750
$ws.Range("M31").Value = -67987.2
$ws.Range("N31").Value = -4340
$ws.Range("H34").Value = 37128.723
$ws.Range("I34").Value = 68282.2
$ws.Range("J34").Value = 3750
$ws.Range("K34").Value = 68282.2
$ws.Range("L34").Value = 3750
$ws.Range("M34").Value = -68080.2
$ws.Range("N34").Value = -4154
$ws.Range("H107").Value = 810
$ws.Range("J107").Value = 569.875
$ws.Range("L107").Value = 569.875
$ws.Range("N107").Value = -4409.875
$ws.Range("H113").Value = 998.3333
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H132").Value = 55559244
$ws.Range("I132").Value = 45458410
$ws.Range("J132").Value = 100002920
$ws.Range("K132").Value = 136375230
$ws.Range("L132").Value = 300008760
$ws.Range("M132").Value = -136372700
$ws.Range("N132").Value = -300013820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 72028.28999999999
$ws.Range("J46").Value = 83983.164
$ws.Range("L46").Value = 251949.492
$ws.Range("N46").Value = -252131.492
$ws.Range("H92").Value = 327.66666
$ws.Range("H122").Value = 3829.1724
$ws.Range("J122").Value = 12961.75
$ws.Range("L122").Value = 116655.75
$ws.Range("N122").Value = -121555.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 355400.28
$ws.Range("I102").Value = 2410.5
$ws.Range("J102").Value = 1202575.8
$ws.Range("K102").Value = 2410.5
$ws.Range("L102").Value = 1202575.8
$ws.Range("M102").Value = -788.5
$ws.Range("N102").Value = -1205819.8
$ws.Range("H113").Value = 2650.9
$ws.Range("J113").Value = 1982.8334
$ws.Range("L113").Value = 1982.8334
$ws.Range("N113").Value = -6322.8334
$ws.Range("H122").Value = 1210.6
$ws.Range("I122").Value = 943.8570999999999
$ws.Range("J122").Value = 1833
$ws.Range("K122").Value = 2831.5713
$ws.Range("L122").Value = 5499
$ws.Range("M122").Value = -381.5712999999996
$ws.Range("N122").Value = -10399

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2268.4614
$ws.Range("I61").Value = 1576.25
$ws.Range("J61").Value = 3376
$ws.Range("K61").Value = 1576.25
$ws.Range("L61").Value = 3376
$ws.Range("M61").Value = -1374.25
$ws.Range("N61").Value = -3780
$ws.Range("H113").Value = 2268.4614
$ws.Range("I113").Value = 1576.25
$ws.Range("J113").Value = 3376
$ws.Range("K113").Value = 1576.25
$ws.Range("L113").Value = 3376
$ws.Range("M113").Value = 593.75
$ws.Range("N113").Value = -7716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1589
$ws.Range("I122").Value = 1220.8
$ws.Range("J122").Value = 1819.125
$ws.Range("K122").Value = 3662.4
$ws.Range("L122").Value = 5457.375
$ws.Range("M122").Value = -1212.4
$ws.Range("N122").Value = -10357.375
$ws.Range("H132").Value = 6022.1113
$ws.Range("I132").Value = 8333.166999999999
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 24999.501
$ws.Range("L132").Value = 4200
$ws.Range("M132").Value = -22469.501
$ws.Range("N132").Value = -9260
